$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row of the A:J data block.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
$lastCol = 10                                                  # columns A..J

for ($i = 1; $i -le 2; $i++) {
    $srcRow = $lastRow
    $newRow = $lastRow + $i

    # Write raw values first -- column A's date advances a day per new row,
    # the rest (B..J) are carried over unchanged from the source row.
    $dateVal = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($newRow, 1).Value = $dateVal + $i

    for ($c = 2; $c -le $lastCol; $c++) {
        $srcVal = $ws.Cells.Item($srcRow, $c).Value2
        $ws.Cells.Item($newRow, $c).Value = $srcVal
    }

    # Then copy the formatting from the source row onto the new row so the
    # same style indices (borders/number format/alignment) are reused.
    $srcRange = $ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, $lastCol))
    $dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
}
